# Auto-generated edit script applying scheduled market-data refresh to Spriggan_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 8483.25
$ws.Range("I9").Value = 9181.727999999999
$ws.Range("J9").Value = 800
$ws.Range("K9").Value = 9181.727999999999
$ws.Range("L9").Value = 800
$ws.Range("M9").Value = -9012.727999999999
$ws.Range("N9").Value = -1138
$ws.Range("H18").Value = 3422.5
$ws.Range("I18").Value = 200
$ws.Range("J18").Value = 4496.6665
$ws.Range("K18").Value = 200
$ws.Range("L18").Value = 4496.6665
$ws.Range("M18").Value = 84
$ws.Range("N18").Value = -5064.6665
$ws.Range("H43").Value = 6370.1665
$ws.Range("I43").Value = 7305.5
$ws.Range("J43").Value = 4499.5
$ws.Range("K43").Value = 7305.5
$ws.Range("L43").Value = 4499.5
$ws.Range("M43").Value = -7236.5
$ws.Range("N43").Value = -4637.5
$ws.Range("H107").Value = 1094.4615
$ws.Range("I107").Value = 1131.9166
$ws.Range("J107").Value = 645
$ws.Range("K107").Value = 1131.9166
$ws.Range("L107").Value = 645
$ws.Range("M107").Value = 788.0834
$ws.Range("N107").Value = -4485
$ws.Range("H116").Value = 5272.478
$ws.Range("J116").Value = 7449.5
$ws.Range("L116").Value = 7449.5
$ws.Range("N116").Value = -14333.5
$ws.Range("H125").Value = 20293948
$ws.Range("I125").Value = 4239035
$ws.Range("J125").Value = 41700500
$ws.Range("K125").Value = 38151315
$ws.Range("L125").Value = 375304500
$ws.Range("M125").Value = -38148855
$ws.Range("N125").Value = -375309420
$ws.Range("H129").Value = 8027.5386
$ws.Range("I129").Value = 1487.3636
$ws.Range("K129").Value = 4462.0908
$ws.Range("M129").Value = 537.9092000000001
$ws.Range("H138").Value = 2570.225
$ws.Range("I138").Value = 1528.1538
$ws.Range("J138").Value = 3071.963
$ws.Range("K138").Value = 4584.4614
$ws.Range("L138").Value = 9215.889000000001
$ws.Range("M138").Value = 555.5385999999999
$ws.Range("N138").Value = -19495.889
$ws.Range("H141").Value = 8514.666999999999
$ws.Range("I141").Value = 8514.666999999999
$ws.Range("K141").Value = 25544.001
$ws.Range("M141").Value = -20364.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4238.7744
$ws.Range("I32").Value = 2374.577
$ws.Range("J32").Value = 13932.6
$ws.Range("K32").Value = 2374.577
$ws.Range("L32").Value = 13932.6
$ws.Range("M32").Value = -2087.577
$ws.Range("N32").Value = -14506.6
$ws.Range("H34").Value = 41900
$ws.Range("I34").Value = 41900
$ws.Range("K34").Value = 41900
$ws.Range("M34").Value = -41629
$ws.Range("H74").Value = 43484240
$ws.Range("I74").Value = 43484240
$ws.Range("K74").Value = 43484240
$ws.Range("M74").Value = -43483366
$ws.Range("H77").Value = 43484240
$ws.Range("I77").Value = 43484240
$ws.Range("K77").Value = 217421200
$ws.Range("M77").Value = -217416832

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 38000
$ws.Range("J30").Value = 38000
$ws.Range("L30").Value = 38000
$ws.Range("N30").Value = -38250
$ws.Range("H31").Value = 3649.5
$ws.Range("I31").Value = 3649.5
$ws.Range("K31").Value = 3649.5
$ws.Range("M31").Value = -3397.5
$ws.Range("H33").Value = 37000
$ws.Range("J33").Value = 37000
$ws.Range("L33").Value = 37000
$ws.Range("N33").Value = -37672
$ws.Range("H37").Value = 10737.5
$ws.Range("I37").Value = 4000
$ws.Range("J37").Value = 17475
$ws.Range("K37").Value = 4000
$ws.Range("L37").Value = 17475
$ws.Range("M37").Value = -3863
$ws.Range("N37").Value = -17749
$ws.Range("H105").Value = 1930.1904
$ws.Range("I105").Value = 1283.375
$ws.Range("K105").Value = 1283.375
$ws.Range("M105").Value = 463.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 202.71428
$ws.Range("I7").Value = 280
$ws.Range("K7").Value = 280
$ws.Range("M7").Value = -167
$ws.Range("H22").Value = 650.25
$ws.Range("I22").Value = 750.5
$ws.Range("J22").Value = 550
$ws.Range("K22").Value = 750.5
$ws.Range("L22").Value = 550
$ws.Range("M22").Value = -400.5
$ws.Range("N22").Value = -1250
$ws.Range("H31").Value = 8551.725
$ws.Range("I31").Value = 5868.4736
$ws.Range("K31").Value = 5868.4736
$ws.Range("M31").Value = -5573.4736
$ws.Range("H34").Value = 8551.725
$ws.Range("I34").Value = 5868.4736
$ws.Range("K34").Value = 5868.4736
$ws.Range("M34").Value = -5666.4736
$ws.Range("H54").Value = 28499.5
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 28499.5
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 28499.5
$ws.Range("M54").Value = ""
$ws.Range("N54").Value = -29815.5
$ws.Range("H92").Value = 95000
$ws.Range("J92").Value = 95000
$ws.Range("L92").Value = 95000
$ws.Range("N92").Value = -99992

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 3366
$ws.Range("J117").Value = 3775.4285
$ws.Range("L117").Value = 11326.2855
$ws.Range("N117").Value = -18210.2855
$ws.Range("H131").Value = 1588.421
$ws.Range("I131").Value = 1165.4166
$ws.Range("K131").Value = 3496.2498
$ws.Range("M131").Value = 1543.7502
$ws.Range("H139").Value = 1530.7084
$ws.Range("I139").Value = 1542.591
$ws.Range("J139").Value = 1400
$ws.Range("K139").Value = 4627.772999999999
$ws.Range("L139").Value = 4200
$ws.Range("M139").Value = 512.2270000000008
$ws.Range("N139").Value = -14480

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").Value = ""
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = ""
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").Value = ""
$ws.Range("H97").Value = 1229.9166
$ws.Range("I97").Value = 917.4706
$ws.Range("J97").Value = 1988.7142
$ws.Range("K97").Value = 917.4706
$ws.Range("L97").Value = 1988.7142
$ws.Range("M97").Value = -421.4706
$ws.Range("N97").Value = -2980.7142
$ws.Range("H102").Value = 3464.5881
$ws.Range("I102").Value = 3464.5881
$ws.Range("K102").Value = 3464.5881
$ws.Range("M102").Value = -1842.5881
$ws.Range("H113").Value = 50634.617
$ws.Range("I113").Value = 69852.2
$ws.Range("J113").Value = 2590.6667
$ws.Range("K113").Value = 69852.2
$ws.Range("L113").Value = 2590.6667
$ws.Range("M113").Value = -67682.2
$ws.Range("N113").Value = -6930.6667
$ws.Range("H122").Value = 56096.39
$ws.Range("I122").Value = 61610.85
$ws.Range("K122").Value = 184832.55
$ws.Range("M122").Value = -182382.55

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 20000.5
$ws.Range("I45").Value = 20000.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 20000.5
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -19593.5
$ws.Range("N45").Value = ""
$ws.Range("H61").Value = 4408.9546
$ws.Range("I61").Value = 4557
$ws.Range("J61").Value = 1300
$ws.Range("K61").Value = 4557
$ws.Range("L61").Value = 1300
$ws.Range("M61").Value = -4355
$ws.Range("N61").Value = -1704
$ws.Range("H113").Value = 4408.9546
$ws.Range("I113").Value = 4557
$ws.Range("J113").Value = 1300
$ws.Range("K113").Value = 4557
$ws.Range("L113").Value = 1300
$ws.Range("M113").Value = -2387
$ws.Range("N113").Value = -5640
$ws.Range("H132").Value = 43641970
$ws.Range("I132").Value = 53338740
$ws.Range("J132").Value = 6497.5
$ws.Range("K132").Value = 160016220
$ws.Range("L132").Value = 19492.5
$ws.Range("M132").Value = -160013690
$ws.Range("N132").Value = -24552.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2213.7778
$ws.Range("I126").Value = 2869.75
$ws.Range("K126").Value = 8609.25
$ws.Range("M126").Value = -6139.25
$ws.Range("H136").Value = 26317306
$ws.Range("I136").Value = 26317306
$ws.Range("K136").Value = 78951918
$ws.Range("M136").Value = -217416832
